$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert the new "login_id_policy_query" column before the old column F
#    (OTP_Value), pushing the old F:W header/data block one column right
#    (to G:X).
# ---------------------------------------------------------------------------
$ws.Columns("F:F").Insert()

# ---------------------------------------------------------------------------
# 2. Insert six new columns (O:T) for the password-change / password-reset
#    query block, pushing the old N:.. block (now shifted to U:AD) further
#    right.
# ---------------------------------------------------------------------------
$ws.Columns("O:T").Insert()

# ---------------------------------------------------------------------------
# 3. Populate the newly inserted cells. The order of assignment matches the
#    order new shared-string entries must be appended in (header then data
#    for column F, then header row then data row for the O:T block), so the
#    rebuilt xl/sharedStrings.xml lines up with the target workbook.
# ---------------------------------------------------------------------------
$ws.Range("F1").Value = "login_id_policy_query"
$ws.Range("F2").Value = "Select P.PARAMTER_VALUE from DC_APPLICATION_PARAM_DETAIL P where P.PARAMETER_NAME = 'USER_NAME_POLICY_DESCRIPTION'"

$ws.Range("O1").Value = "last_pass_change_query"
$ws.Range("P1").Value = "last_tran_pass_change_query"
$ws.Range("Q1").Value = "is_password_change_required_value"
$ws.Range("R1").Value = "is_password_change_required_query"
$ws.Range("S1").Value = "is_password_reset_required_value"
$ws.Range("T1").Value = "is_password_reset_required_query"

$ws.Range("O2").Value = "select P.LAST_PASSWORD_CHANGED from dc_customer_info P where P.CNIC ='{customer_cnic}'"
$ws.Range("P2").Value = "select P.LAST_TRANS_PASSWORD_CHANGED from dc_customer_info P where P.CNIC ='{customer_cnic}'"
$ws.Range("Q2").Value = "0"
$ws.Range("R2").Value = "SELECT P.IS_PASSWORD_CHANGED_REQUIRED from dc_customer_info P where P.CNIC ='{customer_cnic}'"
$ws.Range("S2").Value = "0"
$ws.Range("T2").Value = "SELECT P.IS_PASSWORD_RESET_REQUIRED from dc_customer_info P where P.CNIC ='{customer_cnic}'"

# ---------------------------------------------------------------------------
# 5. Give the six newly inserted columns the same (wide) column width as
#    their neighbour so they render like the rest of the query columns.
# ---------------------------------------------------------------------------
$ws.Range("O1:T2").ColumnWidth = $ws.Columns("N").ColumnWidth

# ---------------------------------------------------------------------------
# 6. Move the view / selection to the new right-hand edge of the sheet.
# ---------------------------------------------------------------------------
$ws.Range("AB1").Select()

# ---------------------------------------------------------------------------
# 7. Re-anchor the existing duplicate-values conditional formatting rule
#    (it stayed pinned to its original I8:I12 range across the column
#    inserts) onto its new location, J8:J12.
# ---------------------------------------------------------------------------
$fc = $ws.Range("I8:I12").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("J8:J12"))
